$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 with the new submission's data
$ws.Range("A2").Value = "Rahav jitthamanyu RA"
$ws.Range("B2").Value = "ra.rahavjitthamanyu@ltimindtree.com"
$ws.Range("C2").Value = "COD"
$ws.Range("D2").Value = 28.2
$ws.Range("E2").Value = 60
$ws.Range("F2").Value = "2025-10-18 | 05:32:22 PM"
$ws.Range("I2").Value = "null"
$ws.Range("J2").Value = "https://admin.ltimindtree.iamneo.ai/result?testId=U2FsdGVkX1%2B1VRj4uLCXzOtOJehrdadk9T3OlVwbQ3TCUKBl8REzy4ZNOseny1IWfhzmyqAXe6HLCrky80lUbmxVvVlPthDW0dAOWDDYMzMrZppBatWZQEReQXY59JqNYladNNWrGIo3f9Y20V2ePA%3D%3D"

# Remove the now-stale rows 3-6 entirely (shifts remaining rows up / shrinks used range)
$ws.Range("A3:J6").EntireRow.Delete()
